$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ------------------------------------------------------------------
# Insert two new data rows in the middle of the sorted table:
#   - row 11: Ca10 / anytx   (Cancer)
#   - row 21: D18  / surgery2 (Environmental)
# Using worksheet row inserts (ListRows.Add ignores its index argument
# in this runtime, so we insert blank rows directly and resize the
# table afterwards).
# ------------------------------------------------------------------
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(21).Insert()

# Fill the two new rows in the same cell-write order the source
# workbook used (A11, A21, B21, D21, B11, D11) so shared-string
# indices line up with the authored file. Category (column C) reuses
# pre-existing shared strings ("Cancer" / "Environmental").
$ws.Range("A11").Value = "Ca10"
$ws.Range("C11").Value = "Cancer"
$ws.Range("A21").Value = "D18"
$ws.Range("C21").Value = "Environmental"
$ws.Range("B21").Value = "surgery2"
$ws.Range("D21").Value = "Derived variable indicating if there has been surgery within 3 months"
$ws.Range("B11").Value = "anytx"
$ws.Range("D11").Value = "Any cancer treatment in the past 3 months"

# ------------------------------------------------------------------
# Append nine new laboratory-variable rows (L2-L10) after the
# existing L1 row, at the bottom of the table.
# ------------------------------------------------------------------
$ws.Range("A75").Value = "L2"
$ws.Range("B75").Value = "ddimer"
$ws.Range("D75").Value = "D-Dimer"

$ws.Range("A76").Value = "L3"
$ws.Range("B76").Value = "fibrinogen"
$ws.Range("D76").Value = "Fibrinogen"

$ws.Range("A77").Value = "L4"
$ws.Range("D77").Value = "PT"
$ws.Range("A78").Value = "L5"
$ws.Range("B77").Value = "pt"
$ws.Range("B78").Value = "aptt"
$ws.Range("D78").Value = "aPTT"

$ws.Range("A79").Value = "L6"
$ws.Range("B79").Value = "hs_trop"
$ws.Range("D79").Value = "High-sensitivity troponin"

$ws.Range("A80").Value = "L7"
$ws.Range("B80").Value = "bnp"
$ws.Range("D80").Value = "BNP"

$ws.Range("A81").Value = "L8"
$ws.Range("B81").Value = "crp"
$ws.Range("D81").Value = "CRP"

$ws.Range("A82").Value = "L9"
$ws.Range("B82").Value = "ldh"
$ws.Range("D82").Value = "LDH"

$ws.Range("A83").Value = "L10"
$ws.Range("B83").Value = "il6"
$ws.Range("D83").Value = "IL-6"

# Column C ("Category") for every new row below row 74 is "Laboratory",
# matching the other L-series rows.
for ($r = 75; $r -le 83; $r++) {
    $ws.Range("C$r").Value = "Laboratory"
}

# ------------------------------------------------------------------
# Grow Table1 to the new extent and refresh the selection to match.
# ------------------------------------------------------------------
$lo.Resize($ws.Range("A1:E83"))
$ws.Range("D83").Select()
